$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I9").Value = "aa"
$ws.Range("J9").Value = "Agree/Accept"
$ws.Range("I14").Value = "%"
$ws.Range("J14").Value = "Uninterpretable"
$ws.Range("I20").Value = "%"
$ws.Range("J20").Value = "Uninterpretable"
$ws.Range("I21").Value = "%"
$ws.Range("J21").Value = "Uninterpretable"
$ws.Range("I22").Value = "%"
$ws.Range("J22").Value = "Uninterpretable"
$ws.Range("I23").Value = "%"
$ws.Range("J23").Value = "Uninterpretable"
$ws.Range("I24").Value = "sd"
$ws.Range("J24").Value = "Statement-non-opinion"
$ws.Range("I31").Value = "sv"
$ws.Range("J31").Value = "Statement-opinion"
$ws.Range("I32").Value = "qy"
$ws.Range("J32").Value = "Yes-No-Question"
$ws.Range("I37").Value = "sd"
$ws.Range("J37").Value = "Statement-non-opinion"
$ws.Range("I39").Value = "ba"
$ws.Range("J39").Value = "Appreciation"
$ws.Range("I40").Value = "sv"
$ws.Range("J40").Value = "Statement-opinion"
$ws.Range("I42").Value = "b"
$ws.Range("J42").Value = "Acknowledge (Backchannel)"
$ws.Range("I47").Value = "b"
$ws.Range("J47").Value = "Acknowledge (Backchannel)"
$ws.Range("I62").Value = "ba"
$ws.Range("J62").Value = "Appreciation"
$ws.Range("I63").Value = "sv"
$ws.Range("J63").Value = "Statement-opinion"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I78").Value = "ba"
$ws.Range("J78").Value = "Appreciation"
$ws.Range("I87").Value = "%"
$ws.Range("J87").Value = "Uninterpretable"
$ws.Range("I91").Value = "sd"
$ws.Range("J91").Value = "Statement-non-opinion"
$ws.Range("I94").Value = "b"
$ws.Range("J94").Value = "Acknowledge (Backchannel)"
$ws.Range("I96").Value = "aa"
$ws.Range("J96").Value = "Agree/Accept"
$ws.Range("I97").Value = "aa"
$ws.Range("J97").Value = "Agree/Accept"
$ws.Range("I101").Value = "sd"
$ws.Range("J101").Value = "Statement-non-opinion"
$ws.Range("I113").Value = "sd"
$ws.Range("J113").Value = "Statement-non-opinion"
$ws.Range("I126").Value = "aa"
$ws.Range("J126").Value = "Agree/Accept"
$ws.Range("I140").Value = "sv"
$ws.Range("J140").Value = "Statement-opinion"
$ws.Range("I160").Value = "sd"
$ws.Range("J160").Value = "Statement-non-opinion"
$ws.Range("I162").Value = "b"
$ws.Range("J162").Value = "Acknowledge (Backchannel)"
$ws.Range("I165").Value = "aa"
$ws.Range("J165").Value = "Agree/Accept"
$ws.Range("I166").Value = "aa"
$ws.Range("J166").Value = "Agree/Accept"
$ws.Range("I170").Value = "sd"
$ws.Range("J170").Value = "Statement-non-opinion"
$ws.Range("I179").Value = "b"
$ws.Range("J179").Value = "Acknowledge (Backchannel)"
$ws.Range("I184").Value = "sd"
$ws.Range("J184").Value = "Statement-non-opinion"
$ws.Range("I189").Value = "sd"
$ws.Range("J189").Value = "Statement-non-opinion"
$ws.Range("I190").Value = "sd"
$ws.Range("J190").Value = "Statement-non-opinion"
$ws.Range("I196").Value = "qy"
$ws.Range("J196").Value = "Yes-No-Question"
$ws.Range("I207").Value = "b"
$ws.Range("J207").Value = "Acknowledge (Backchannel)"
$ws.Range("I208").Value = "sv"
$ws.Range("J208").Value = "Statement-opinion"
$ws.Range("I226").Value = "sd"
$ws.Range("J226").Value = "Statement-non-opinion"
$ws.Range("I236").Value = "ba"
$ws.Range("J236").Value = "Appreciation"
$ws.Range("I246").Value = "aa"
$ws.Range("J246").Value = "Agree/Accept"
$ws.Range("I251").Value = "sv"
$ws.Range("J251").Value = "Statement-opinion"
$ws.Range("I252").Value = "sd"
$ws.Range("J252").Value = "Statement-non-opinion"
$ws.Range("I263").Value = "aa"
$ws.Range("J263").Value = "Agree/Accept"
$ws.Range("I265").Value = "aa"
$ws.Range("J265").Value = "Agree/Accept"
$ws.Range("I266").Value = "sd"
$ws.Range("J266").Value = "Statement-non-opinion"
$ws.Range("I270").Value = "aa"
$ws.Range("J270").Value = "Agree/Accept"
$ws.Range("I274").Value = "sd"
$ws.Range("J274").Value = "Statement-non-opinion"
$ws.Range("I289").Value = "ba"
$ws.Range("J289").Value = "Appreciation"
$ws.Range("I297").Value = "sd"
$ws.Range("J297").Value = "Statement-non-opinion"
$ws.Range("I301").Value = "b"
$ws.Range("J301").Value = "Acknowledge (Backchannel)"
$ws.Range("I302").Value = "sd"
$ws.Range("J302").Value = "Statement-non-opinion"
$ws.Range("I304").Value = "sd"
$ws.Range("J304").Value = "Statement-non-opinion"
$ws.Range("I310").Value = "aa"
$ws.Range("J310").Value = "Agree/Accept"
$ws.Range("I316").Value = "aa"
$ws.Range("J316").Value = "Agree/Accept"
$ws.Range("I318").Value = "sv"
$ws.Range("J318").Value = "Statement-opinion"
$ws.Range("I319").Value = "%"
$ws.Range("J319").Value = "Uninterpretable"
$ws.Range("I320").Value = "sd"
$ws.Range("J320").Value = "Statement-non-opinion"
$ws.Range("I321").Value = "sd"
$ws.Range("J321").Value = "Statement-non-opinion"
$ws.Range("I329").Value = "sd"
$ws.Range("J329").Value = "Statement-non-opinion"
$ws.Range("I330").Value = "aa"
$ws.Range("J330").Value = "Agree/Accept"
$ws.Range("I347").Value = "sd"
$ws.Range("J347").Value = "Statement-non-opinion"
$ws.Range("I348").Value = "sv"
$ws.Range("J348").Value = "Statement-opinion"
$ws.Range("I354").Value = "sd"
$ws.Range("J354").Value = "Statement-non-opinion"
$ws.Range("I356").Value = "sd"
$ws.Range("J356").Value = "Statement-non-opinion"
$ws.Range("I357").Value = "aa"
$ws.Range("J357").Value = "Agree/Accept"
$ws.Range("I366").Value = "ba"
$ws.Range("J366").Value = "Appreciation"
$ws.Range("I396").Value = "ba"
$ws.Range("J396").Value = "Appreciation"
$ws.Range("I399").Value = "sd"
$ws.Range("J399").Value = "Statement-non-opinion"
$ws.Range("I413").Value = "aa"
$ws.Range("J413").Value = "Agree/Accept"
$ws.Range("I415").Value = "sd"
$ws.Range("J415").Value = "Statement-non-opinion"
$ws.Range("I420").Value = "%"
$ws.Range("J420").Value = "Uninterpretable"
$ws.Range("I421").Value = "%"
$ws.Range("J421").Value = "Uninterpretable"
$ws.Range("I439").Value = "sd"
$ws.Range("J439").Value = "Statement-non-opinion"
$ws.Range("I451").Value = "sv"
$ws.Range("J451").Value = "Statement-opinion"
$ws.Range("I462").Value = "ba"
$ws.Range("J462").Value = "Appreciation"
$ws.Range("I469").Value = "sd"
$ws.Range("J469").Value = "Statement-non-opinion"
$ws.Range("I471").Value = "aa"
$ws.Range("J471").Value = "Agree/Accept"
$ws.Range("I475").Value = "b"
$ws.Range("J475").Value = "Acknowledge (Backchannel)"
